# "add function is working"
#
# On the "vendors" sheet, a new vendor row ("test vendor") was added above
# the existing "Kate's Car and Co." entry (to confirm that the "add vendor"
# function works), and the "Kate's Car and Co." entry itself was retyped
# without the apostrophe / trailing line breaks ("Kates Car and Co.").
# The author also briefly visited the "Kates Car and Co." sheet (leaving its
# selection on C14) before returning focus to the "vendors" sheet.

$wb = $excel.ActiveWorkbook

$wsVendors = $wb.Worksheets.Item("vendors")
$wsKates   = $wb.Worksheets.Item("Kates Car and Co.")

# Insert a new row 2 on "vendors" and fill it in, pushing the old rows down.
$wsVendors.Rows("2:2").Insert()

# Row 3 now holds the old "Kate's Car and Co." text - retype it cleanly.
$wsVendors.Range("A3").Value = "Kates Car and Co."

# Row 2 is the freshly-inserted test row.
$wsVendors.Range("A2").Value = "test vendor"

# Match the row heights Excel settled on after the edit/autofit.
$wsVendors.Rows("1:2").RowHeight = 15.75
$wsVendors.Rows("3:5").RowHeight = 12.5

# Visit "Kates Car and Co." (selecting C14), then return to "vendors".
$wsKates.Range("C14").Select()
$wsVendors.Select()
$wsVendors.Range("B6").Select()
